$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ParticipantsTab"
$ws.Range("A2").Select()
